# Automatic update of files.
# The source rows (2-13, excluding unchanged row 11) get their data replaced
# by another row's original data (a row-content permutation), as captured
# in the xml diff. Build the target state per row and write it back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; A=111471685; B=77515; D="NT"; E=6425; F="Garnlav"; G="Alectoria sarmentosa"; H="(Ach.) Ach."; Q=554595.0694405095; R=7003142.694495555; S=25; Z="15:49"; AB="15:49"; AC=$null },
    @{ Row=3; A=111471083; B=96348; D="VU"; E=220787; F="Knärot"; G="Goodyera repens"; H="(L.) R. Br."; Q=554499.1143642976; R=7003141.52872613; S=25; Z="15:31"; AB="15:31"; AC=$null },
    @{ Row=4; A=111470448; B=96348; D="VU"; E=220787; F="Knärot"; G="Goodyera repens"; H="(L.) R. Br."; Q=554488.5866359913; R=7003175.257923778; S=22; Z="14:59"; AB="14:59"; AC=$null },
    @{ Row=5; A=111470636; B=94134; D="NT"; E=53; F="Vedtrappmossa"; G="Crossocalyx hellerianus"; H="(Nees ex Lindenb.) Meyl."; Q=554457.9939421143; R=7003163.892755959; S=25; Z="14:41"; AB="14:41"; AC=$null },
    @{ Row=6; A=111470743; B=78611; D="LC"; E=6463; F="Bårdlav"; G="Nephroma parile"; H="(Ach.) Ach."; Q=554457.9939421143; R=7003163.892755959; S=25; Z="14:41"; AB="14:41"; AC=$null },
    @{ Row=7; A=111471797; B=77515; D="NT"; E=6425; F="Garnlav"; G="Alectoria sarmentosa"; H="(Ach.) Ach."; Q=554597.2688619854; R=7003280.616068945; S=25; Z="15:49"; AB="15:49"; AC="På tall" },
    @{ Row=8; A=111470792; B=96348; D="VU"; E=220787; F="Knärot"; G="Goodyera repens"; H="(L.) R. Br."; Q=554440.9784625648; R=7003152.756292564; S=25; Z="15:19"; AB="15:19"; AC=$null },
    @{ Row=9; A=111470486; B=78578; D="NT"; E=6458; F="Lunglav"; G="Lobaria pulmonaria"; H="(L.) Hoffm."; Q=554488.5866359913; R=7003175.257923778; S=22; Z="14:41"; AB="14:41"; AC=$null },
    @{ Row=10; A=111470101; B=77515; D="NT"; E=6425; F="Garnlav"; G="Alectoria sarmentosa"; H="(Ach.) Ach."; Q=554474.9281677724; R=7003314.266989549; S=25; Z="00:00"; AB="00:00"; AC="På tall" },
    @{ Row=12; A=111469986; B=77515; D="NT"; E=6425; F="Garnlav"; G="Alectoria sarmentosa"; H="(Ach.) Ach."; Q=554489.6113782075; R=7003329.432399829; S=25; Z="00:00"; AB="00:00"; AC="Rikligt på tall" },
    @{ Row=13; A=111470685; B=77267; D="NT"; E=6446; F="Kolflarnlav"; G="Carbonicola anthracophila"; H="(Nyl.) Bendiksby & Timdal"; Q=554457.9939421143; R=7003163.892755959; S=25; Z="14:41"; AB="14:41"; AC=$null }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("Z$row").Value = $r.Z
    $ws.Range("AB$row").Value = $r.AB
    if ($null -eq $r.AC) {
        $ws.Range("AC$row").ClearContents()
    } else {
        $ws.Range("AC$row").Value = $r.AC
    }
}
